$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header row (row 1): add five new date columns (C:G) between the existing
# "18_12_2023" (B) and "07_02_2024"/"11_02_2024"/"18_02_2024" columns, and
# shift the two previously-existing data columns out to H and I.
# ---------------------------------------------------------------------------
$ws.Range("B1").Value = "18_12_2023"
$ws.Range("C1").Value = "05_01_2024"
$ws.Range("D1").Value = "15_01_2024"
$ws.Range("E1").Value = "21_01_2024"
$ws.Range("F1").Value = "28_01_2024"
$ws.Range("G1").Value = "07_02_2024"
$ws.Range("H1").Value = "11_02_2024"
$ws.Range("I1").Value = "18_02_2024"

# ---------------------------------------------------------------------------
# Data rows 2-6: columns B-G are brand-new weekly counts, H keeps the old
# "B" column value and I keeps the old "C" column value.
# ---------------------------------------------------------------------------

# Row 2 - Alejandro
$ws.Range("B2").Value = 1826
$ws.Range("C2").Value = 2293
$ws.Range("D2").Value = 2525
$ws.Range("E2").Value = 2638
$ws.Range("F2").Value = 2824
$ws.Range("G2").Value = 3063
$ws.Range("H2").Value = 3215
$ws.Range("I2").Value = 3215

# Row 3 - Camila
$ws.Range("B3").Value = 1269
$ws.Range("C3").Value = 1716
$ws.Range("D3").Value = 1910
$ws.Range("E3").Value = 2032
$ws.Range("F3").Value = 2097
$ws.Range("G3").Value = 2314
$ws.Range("H3").Value = 2385
$ws.Range("I3").Value = 2640

# Row 4 - Betty
$ws.Range("B4").Value = 3019
$ws.Range("C4").Value = 3196
$ws.Range("D4").Value = 3373
$ws.Range("E4").Value = 3446
$ws.Range("F4").Value = 3537
$ws.Range("G4").Value = 3618
$ws.Range("H4").Value = 3682
$ws.Range("I4").Value = 3769

# Row 5 - Felipe
$ws.Range("B5").Value = 6672
$ws.Range("C5").Value = 6688
$ws.Range("D5").Value = 6704
$ws.Range("E5").Value = 6784
$ws.Range("F5").Value = 6858
$ws.Range("G5").Value = 6927
$ws.Range("H5").Value = 6931
$ws.Range("I5").Value = 7034

# Row 6 - Constanza
$ws.Range("B6").Value = 0
$ws.Range("C6").Value = 0
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 33

# ---------------------------------------------------------------------------
# Formatting: the underlined-font style that used to sit on C3 now belongs to
# I3 (the cell holding the value that was formerly in C3); make sure C3 goes
# back to the default (no underline).
# ---------------------------------------------------------------------------
$ws.Range("C3").Font.Underline = -4142
$ws.Range("I3").Font.Underline = 2

# ---------------------------------------------------------------------------
# Column widths: B and C get explicit widths (closest achievable via the
# ColumnWidth property, which snaps to pixel-sized increments).
# ---------------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 12.35
$ws.Columns.Item(3).ColumnWidth = 13.0

# ---------------------------------------------------------------------------
# Selection: active cell moves to I3.
# ---------------------------------------------------------------------------
$ws.Range("I3").Select()
